$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the brand-new "PEAS 8PUZZLE" section at the very start of the
#    document, before the existing "PEAS CỜ CARO" title paragraph.
# ---------------------------------------------------------------------------
$insertXml = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>PEAS 8PUZZLE</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>1. Performance </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Mục tiêu: Đưa trạng thái ban đầu về trạng thái đích (thường là dãy số từ 1-8 theo thứ tự, ô trống ở cuối, hoặc hình ảnh nào đó)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Hạn chế số bước để thắng (tìm chiến thắng nhanh nhất).</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Rút ngắn thời gian giải</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>2. Environment </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Không gian</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Tất cả các hoán vị hợp lệ của 8 quân và 1 ô trống trên bàn 3x3.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Trạng thái ban đầu: </w:t>
  </w:r>
  <w:r>
    <w:t>Một cấu hình bất kì của 8 ô và 1 ô trống</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Trạng thái:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Cấu hình mong muốn (Dãy số 12345678_ hoặc hình ảnh hoàn chỉnh</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t xml:space="preserve"> 3. Actuators</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Tác nhân có thể di chuyển ô trống theo 4 hướng (Left, Right, Up, Down) nếu hợp lệ</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t xml:space="preserve"> 4. Sensors </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Quan sát toàn bộ </w:t>
  </w:r>
  <w:r>
    <w:t>bảng</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> hiện tại:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Vị trí hiện tại của các ô (bảng 3x3)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Vị trí ô trống</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></w:rPr>
    <w:sym w:font="Wingdings" w:char="F0E8"/>
  </w:r>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Tóm gọn PEAS của</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t xml:space="preserve"> cờ</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t xml:space="preserve"> Caro:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>P:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Giải xong với số bước ít nhất / chi phí thấp nhất.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>E:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Bảng 3x3 với 8 quân số và 1 ô trống.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>A:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Di chuyển ô trống (Up, Down, Left, Right)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>S:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:t>Quan sát toàn bộ trạng thái bảng (vị trí của các quân số + ô trống).</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$start = $d.Range(0, 0)
$start.InsertXML($insertXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) The original "PEAS CỜ CARO" title paragraph (now pushed down after the
#    newly-inserted section) gains a <w:lastRenderedPageBreak/> marker,
#    since it now starts a new page.
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("PEAS CỜ CARO", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $caroTitle = $d.Content
    $caroRange = $d.Range($d.Content.Find.Parent.Start, $d.Content.Find.Parent.End)
}

# Locate paragraph whose text is exactly "PEAS CỜ CARO" and rebuild its single
# run so it includes the lastRenderedPageBreak element.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "PEAS CỜ CARO") {
        $r = $p.Range
        $textRange = $d.Range($r.Start, $r.End - 1)
        $xml = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:r>
<w:rPr><w:b/><w:bCs/></w:rPr>
<w:lastRenderedPageBreak/>
<w:t>PEAS CỜ CARO</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
        $textRange.InsertXML($xml) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Merge the two runs "PEAS " + "TRÒ CHƠI 8 QUÂN HẬU" of the second title
#    into a single run "PEAS TRÒ CHƠI 8 QUÂN HẬU" (keeping the existing
#    lastRenderedPageBreak on that run).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "PEAS TRÒ CHƠI 8 QUÂN HẬU") {
        $r = $p.Range
        $textRange = $d.Range($r.Start, $r.End - 1)
        $xml2 = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:body>
<w:p>
<w:r>
<w:rPr><w:b/><w:bCs/></w:rPr>
<w:lastRenderedPageBreak/>
<w:t>PEAS TRÒ CHƠI 8 QUÂN HẬU</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
        $textRange.InsertXML($xml2) | Out-Null
        break
    }
}
